# import-clients.xlsx template — re-order the "Client 2" / "Address" header
# columns: the "address_1" header (previously column I) moves to column E,
# pushing the "first_name_2 / last_name_2 / email_2 / phone_2" headers
# (previously E:H) one column to the right, into F:I. Columns J onward
# ("address_2", "city", "state", "zip", "home_type", "notes") are untouched.
#
# Doing this as a column Cut + "Insert Cut Cells" (exactly the gesture a
# user performs in the Excel UI: select column I, Cut, select column E,
# Insert Cut Cells) keeps the sheet's total column count the same and
# carries the "address_1" header's cell formatting along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("I:I").Cut()
$ws.Columns("E:E").Insert()

# Match the saved selection/active-cell state from the edited workbook.
$ws.Range("K13").Select()
